# "results code, bad written"
# Re-worked the VRP node-priority results on the "nodes" sheet: several rows
# had the wrong "priority" (column E) value recorded, so fix them up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nodes")
$ws.Activate()

# Column E = "priority" -- correct the mis-recorded results
$ws.Range("E7").Value  = 0   # was 1
$ws.Range("E10").Value = 0   # was 1
$ws.Range("E11").Value = 0   # was 2
$ws.Range("E15").Value = 0   # was 1
$ws.Range("E16").Value = 2   # was 0
$ws.Range("E17").Value = 2   # was 1
$ws.Range("E18").Value = 2   # was 0

# Leave the view scrolled/selected where the edit finished
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H17").Select()
